$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the two new values (commit message: "7 and 9 edited from D drive" -
# rows where column A is 7 and 9, i.e. rows 8 and 10)
$ws.Range("B8").Value = 17
$ws.Range("B10").Value = 19

# Move the selection/active cell to B11 (was B13)
$ws.Range("B11").Select()
